$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-18 07:18:32"
$ws.Range("H2").Value = "'73%"
$ws.Range("N2").Value = "-1.9 °C 6:50 TU"
$ws.Range("E3").Value = "2026-02-18 07:18:34"
$ws.Range("H3").Value = "'95%"
$ws.Range("M3").Value = "-2.6 °C 6:57 TU"
$ws.Range("E4").Value = "2026-02-18 07:18:36"
$ws.Range("O4").Value = "6.7 °C"
$ws.Range("E5").Value = "2026-02-18 07:18:39"
$ws.Range("H5").Value = "'85%"
$ws.Range("M5").Value = "0.8 °C 6:58 TU"
$ws.Range("O5").Value = "-1.6 °C"
$ws.Range("E6").Value = "2026-02-18 07:18:41"
$ws.Range("H6").Value = "'97%"
$ws.Range("J6").Value = "1017.9 hPa"
$ws.Range("N6").Value = "6.0 °C 6:59 TU"
$ws.Range("O6").Value = "7.7 °C"
$ws.Range("E7").Value = "2026-02-18 07:18:44"
$ws.Range("J7").Value = "1018.1 hPa"
$ws.Range("O7").Value = "12.2 °C"
$ws.Range("E8").Value = "2026-02-18 07:18:46"
$ws.Range("J8").Value = "1018.2 hPa"
$ws.Range("L8").Value = "16.2 km/h - 43º 0:12 TU"
$ws.Range("E9").Value = "2026-02-18 07:18:48"
$ws.Range("O9").Value = "4.8 °C"
$ws.Range("E10").Value = "2026-02-18 07:18:51"
$ws.Range("K10").Value = "0.0 MJ/m2"
$ws.Range("N10").Value = "4.1 °C 6:36 TU"
$ws.Range("O10").Value = "6.4 °C"
$ws.Range("E11").Value = "2026-02-18 07:18:53"
$ws.Range("N11").Value = "-0.3 °C 6:48 TU"
$ws.Range("O11").Value = "1.6 °C"
$ws.Range("E12").Value = "2026-02-18 07:18:55"
$ws.Range("O12").Value = "5.7 °C"
$ws.Range("E13").Value = "2026-02-18 07:18:58"
$ws.Range("N13").Value = "-4.4 °C 6:52 TU"
$ws.Range("O13").Value = "-2.5 °C"
$ws.Range("E14").Value = "2026-02-18 07:19:00"
$ws.Range("H14").Value = "'100%"
$ws.Range("O14").Value = "9.8 °C"
$ws.Range("E15").Value = "2026-02-18 07:19:03"
$ws.Range("H15").Value = "'95%"
$ws.Range("O15").Value = "5.3 °C"
$ws.Range("E16").Value = "2026-02-18 07:19:05"
$ws.Range("E17").Value = "2026-02-18 07:19:07"
$ws.Range("H17").Value = "'85%"
$ws.Range("K17").Value = "0.0 MJ/m2"
$ws.Range("E18").Value = "2026-02-18 07:19:09"
$ws.Range("J18").Value = "1018.2 hPa"
$ws.Range("L18").Value = "7.6 km/h - 254º 6:31 TU"
$ws.Range("O18").Value = "7.1 °C"
$ws.Range("E19").Value = "2026-02-18 07:19:12"
$ws.Range("N19").Value = "5.1 °C 6:55 TU"
$ws.Range("O19").Value = "5.4 °C"
$ws.Range("E20").Value = "2026-02-18 07:19:14"
$ws.Range("H20").Value = "'76%"
$ws.Range("E21").Value = "2026-02-18 07:19:17"
$ws.Range("N21").Value = "0.1 °C 6:57 TU"
$ws.Range("O21").Value = "1.6 °C"
$ws.Range("E22").Value = "2026-02-18 07:19:19"
$ws.Range("E23").Value = "2026-02-18 07:19:21"
$ws.Range("E24").Value = "2026-02-18 07:19:24"
$ws.Range("O24").Value = "4.8 °C"
$ws.Range("E25").Value = "2026-02-18 07:19:26"
$ws.Range("M25").Value = "2.3 °C 6:40 TU"
$ws.Range("O25").Value = "-0.1 °C"
$ws.Range("E26").Value = "2026-02-18 07:19:28"
$ws.Range("E27").Value = "2026-02-18 07:19:30"
$ws.Range("H27").Value = "'49%"
$ws.Range("K27").Value = "0.0 MJ/m2"
$ws.Range("E28").Value = "2026-02-18 07:19:33"
$ws.Range("J28").Value = "1018.6 hPa"
$ws.Range("O28").Value = "4.5 °C"
$ws.Range("E29").Value = "2026-02-18 07:19:35"
$ws.Range("O29").Value = "9.2 °C"
$ws.Range("E30").Value = "2026-02-18 07:19:37"
$ws.Range("J30").Value = "1018.2 hPa"
$ws.Range("N30").Value = "5.0 °C 6:39 TU"
$ws.Range("O30").Value = "5.8 °C"
$ws.Range("E31").Value = "2026-02-18 07:19:40"
$ws.Range("H31").Value = "'79%"
$ws.Range("J31").Value = "1016.8 hPa"
$ws.Range("N31").Value = "9.5 °C 6:59 TU"
$ws.Range("O31").Value = "10.4 °C"
$ws.Range("E32").Value = "2026-02-18 07:19:42"
$ws.Range("E33").Value = "2026-02-18 07:19:44"
$ws.Range("H33").Value = "'85%"
$ws.Range("J33").Value = "1021.1 hPa"
$ws.Range("E34").Value = "2026-02-18 07:19:47"
$ws.Range("H34").Value = "'51%"
$ws.Range("L34").Value = "23.8 km/h - 21º 6:52 TU"
$ws.Range("O34").Value = "0.5 °C"
$ws.Range("E35").Value = "2026-02-18 07:19:49"
$ws.Range("H35").Value = "'83%"
$ws.Range("N35").Value = "3.8 °C 6:53 TU"
$ws.Range("O35").Value = "6.8 °C"
$ws.Range("E36").Value = "2026-02-18 07:19:51"
$ws.Range("J36").Value = "1018.0 hPa"
$ws.Range("O36").Value = "8.1 °C"
$ws.Range("E37").Value = "2026-02-18 07:19:54"
$ws.Range("N37").Value = "-0.3 °C 6:47 TU"
$ws.Range("O37").Value = "1.1 °C"
$ws.Range("E38").Value = "2026-02-18 07:19:56"
$ws.Range("N38").Value = "5.0 °C 6:51 TU"
$ws.Range("O38").Value = "8.3 °C"
$ws.Range("E39").Value = "2026-02-18 07:19:58"
$ws.Range("O39").Value = "0.5 °C"
$ws.Range("E40").Value = "2026-02-18 07:20:01"
$ws.Range("J40").Value = "1021.5 hPa"
$ws.Range("N40").Value = "-0.9 °C 6:31 TU"
$ws.Range("O40").Value = "0.5 °C"
$ws.Range("E41").Value = "2026-02-18 07:20:03"
$ws.Range("J41").Value = "1017.8 hPa"
$ws.Range("K41").Value = "0.0 MJ/m2"
$ws.Range("O41").Value = "8.1 °C"
$ws.Range("E42").Value = "2026-02-18 07:20:05"
$ws.Range("O42").Value = "7.8 °C"
$ws.Range("E43").Value = "2026-02-18 07:20:08"
$ws.Range("H43").Value = "'99%"
$ws.Range("L43").Value = "11.2 km/h - 188º 6:54 TU"
$ws.Range("N43").Value = "4.5 °C 6:42 TU"
$ws.Range("O43").Value = "6.7 °C"
$ws.Range("E44").Value = "2026-02-18 07:20:10"
$ws.Range("H44").Value = "'63%"
$ws.Range("M44").Value = "-0.7 °C 6:59 TU"
$ws.Range("O44").Value = "-3.4 °C"
$ws.Range("E45").Value = "2026-02-18 07:20:12"
$ws.Range("J45").Value = "1020.2 hPa"
$ws.Range("M45").Value = "2.2 °C 6:48 TU"
$ws.Range("E46").Value = "2026-02-18 07:20:14"
$ws.Range("J46").Value = "1018.8 hPa"
